$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (KELAS) for all data rows (2-14) to the new value "asdasdsa    "
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = "asdasdsa    "
}

# Update the selected cell / active cell on the sheet
$ws.Range("D16").Select()
